$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column D
$ws.Range("D1").Value = "tier"

# Set the formula for D3:D32 first as one shared-formula block (relative ref
# adjusts per row), then set D2 individually so it stays its own formula
# instance - this reproduces the exact shared-formula grouping Excel wrote
# (D2 standalone, D3:D32 sharing si=0).
$ws.Range("D3:D32").Formula = "=IF(A3<10, 1, IF(A3<100, 2, IF(A3<10000, 3, 4)))"
$ws.Range("D2").Formula = "=IF(A2<10, 1, IF(A2<100, 2, IF(A2<10000, 3, 4)))"

# Update the selection to match the target state
$ws.Range("D2:D32").Select()
$excel.ActiveCell = $ws.Range("D32")

# Adjust the window position (xWindow changed from 4000 to 860)
$excel.Windows.Item(1).Left = 860
